$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update the Date value (row 8, column B) ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# --- Sheet "Include #0": update the System URI value (row 4, column B) ---
$wsInc0 = $wb.Worksheets.Item("Include #0")
$wsInc0.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R44-QualificationPAC"

# --- Sheet "Include #1": update the System URI value (row 4, column B) ---
$wsInc1 = $wb.Worksheets.Item("Include #1")
$wsInc1.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R01-EnsembleSavoirFaire-CISIS"
